# Generate Report for Handoff
# Replaces the old GUID-named source file ("f9f2cc28-...") with the new one
# ("3dd44432-...") across the Overview / zh-cn / de-de sheets, refreshes the
# handoff timestamps + xlf hash, and clears out the now-stale "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" values since a
# fresh handoff has no handback yet (Has metadata flips True).

$wb = $excel.ActiveWorkbook

$oldGuid = "f9f2cc28-e7d9-4043-b141-bf9d022e8f98"
$newGuid = "3dd44432-54e2-4608-8739-fd7176d5960d"
$newHash = "465fc9ac2634be7e3a8e5d38871625454b0b73e7"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-11-09 01:12:49"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-11-09 01:12:35"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("O2").Value = "True"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-11-09 01:12:49"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("O2").Value = "True"
